$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1").Value = "VOLA 3"
$d1 = $ws.Range("D1").Value()
Write-Host "D1: $d1"
